$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.632.21'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '1.591.31'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.95'
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  -2.12%  '
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = '1.813.71'
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").Value = '1.591.94'
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("E15").Value = '  -2.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.68'
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("D17").Value = '26.636.50'
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.61'
$ws.Range("E19").Value = '  -3.33%  '
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.73'
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("E22").Value = '  -2.40%  '
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("E24").Value = '  -0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.86'
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.26'
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  -3.08%  '
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.663'
$ws.Range("E33").Value = '  +22.05%  '
$ws.Range("E34").Value = '  -2.30%  '
$ws.Range("D35").Value = '1.320.44'
$ws.Range("E35").Value = '  -1.08%  '
$ws.Range("E36").Value = '  -4.05%  '
$ws.Range("E37").Value = '  -2.36%  '
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.831'
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("E40").Value = '  -0.24%  '
$ws.Range("E41").Value = '  +3.50%  '
$ws.Range("E42").Value = '  -1.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.11'
$ws.Range("E44").Value = '  -1.84%  '
$ws.Range("D45").Value = '1.726.68'
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.05'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.837'
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.53'
$ws.Range("E51").Value = '  +0.18%  '
